$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 changes from text label "Unnamed: 0" to numeric value 1
$ws.Range("B1").Value = 1

$cols = @("B","C","D","E","F","G","H","I","J","K","L")

$data = @{}
$data[2] = @(2397, 0, 149, 504, 13, 3, 10.64, 191.11, 2, 339, 244.4)
$data[3] = @(6, 0, 329, 334, 4, 0, 1.2, 65.85, 3, 0, 0)
$data[4] = @(269, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[5] = @(24, 0, 155, 155, 0, 0, 0, 71.43, 2, 0, 0)
$data[6] = @(268, 0, 11, 11, 0, 0, 0, 63.64, 0, 0, 0)
$data[7] = @(272, 0, 10, 10, 0, 0, 0, 70, 0, 0, 0)
$data[8] = @(55, 0, 153, 155, 2, 0, 1.29, 87.16, 0, 0, 0)
$data[9] = @(53, 0, 24, 24, 0, 0, 0, 87.5, 0, 0, 0)
$data[10] = @(810, 0, 24, 25, 1, 0, 9.09, 177.78, 1, 0, 0)
$data[11] = @(2458, 0, 159, 164, 5, 1, 17.24, 254.14, 6, 0, 0)
$data[12] = @(1678, 0, 14, 17, 3, 0, 50, 190, 1, 0, 0)
$data[13] = @(956, 1, 154, 178, 23, 0, 71.22, 51.7, 2, 0, 0)
$data[14] = @(23, 0, 125, 126, 1, 0, 0.79, 76.86, 1, 0, 0)
$data[15] = @(270, 0, 4, 4, 0, 0, 0, 75, 0, 0, 0)
$data[16] = @(1114, 0, 217, 270, 11, 0, 9.07, 133.7, 2, 40, 19.05)
$data[17] = @(409, 0, 11, 11, 0, 0, 0, 54.55, 0, 0, 0)
$data[18] = @(295, 0, 298, 314, 5, 1, 1.6, 53.38, 4, 10, 3.18)
$data[19] = @(139, 0, 87, 88, 1, 0, 1.27, 104.7, 1, 0, 0)
$data[20] = @(655, 0, 85, 101, 14, 2, 15.22, 73.44, 2, 0, 0)
$data[21] = @(1473, 0, 32, 33, 0, 1, 0, 301.82, 1, 0, 0)
$data[22] = @(321, 0, 1, 1, 0, 0, 0, 0, 0, 0, 0)
$data[23] = @(14, 0, 7, 7, 0, 0, 0, 100, 0, 0, 0)
$data[24] = @(1284, 0, 19, 19, 0, 0, 0, 240.91, 0, 0, 0)
$data[25] = @(403, 0, 30, 33, 1, 0, 3.03, 53.33, 0, 2, 6.06)
$data[26] = @(354, 0, 12, 15, 1, 1, 7.14, 83.33, 0, 0, 0)
$data[27] = @(325, 0, 41, 42, 0, 1, 0, 24.39, 1, 0, 0)

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $rowValues[$i]
    }
}
